$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (unchanged labels, just re-pointed) -----------------------
$ws.Cells.Item(1, 1).Value = "no"
$ws.Cells.Item(1, 2).Value = "name"
$ws.Cells.Item(1, 3).Value = "notes"

# --- PENJASORKES split into BOYS / GIRLS -----------------------------------
# Written first (ahead of the rest of the alphabetised list) so the two new
# strings land in the shared-string table right after the existing ones, the
# same order the source workbook append used.
$ws.Cells.Item(18, 2).Value = "PENJASORKES BOYS"
$ws.Cells.Item(19, 2).Value = "PENJASORKES GIRLS"

# --- Full (re-alphabetised) course list for B2:B30 -------------------------
$names = @{
    2  = "AGAMA"
    3  = "ART"
    4  = "BAHASA INDONESIA"
    5  = "BAHASA INGGRIS"
    6  = "BIOLOGI"
    7  = "BIOLOGI INTERNATIONAL"
    8  = "DEBAT"
    9  = "EKONOMI"
    10 = "FISIKA"
    11 = "FISIKA INTERNATIONAL"
    12 = "GEOGRAFI"
    13 = "KIMIA"
    14 = "LITERATUR"
    15 = "LITERATUR INDONESIA"
    16 = "MATEMATIKA"
    17 = "MATEMATIKASS"
    20 = "SEJARAH"
    21 = "SENI BUDAYA"
    22 = "SERVICE ELECTIVE"
    23 = "SOSIOLOGI"
    24 = "TIK"
    25 = "UN EKONOMI"
    26 = "UN MATEMATIKASS"
    27 = "UN MATEMATIKA"
    28 = "WALI KELAS"
    29 = "ENGLISH TUTORING"
    30 = "SPANYOL"
}

for ($row = 2; $row -le 30; $row++) {
    if ($names.ContainsKey($row)) {
        $ws.Cells.Item($row, 2).Value = $names[$row]
    }
}

# --- Leftover formatting from the editing session ---------------------------
# (closest reproducible width to the source file's 19.42578125 given this
# host's column-width quantisation)
$ws.Columns.Item(6).ColumnWidth = 18.6
$ws.Range("F1:G1048576").Select() | Out-Null

$wb.Save()
